$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(133).Insert()

$ws.Cells.Item(133, 1).Value = 10
$ws.Cells.Item(133, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(133, 3).Value = "La Araucanía"
$ws.Cells.Item(133, 4).Value = 44855
$ws.Cells.Item(133, 5).Value = 9
$ws.Cells.Item(133, 6).Value = "Fruta"
$ws.Cells.Item(133, 7).Value = 100107
$ws.Cells.Item(133, 8).Value = "Otros"
$ws.Cells.Item(133, 9).Value = 100107002
$ws.Cells.Item(133, 10).Value = "Chirimoya"
$ws.Cells.Item(133, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(133, 12).Value = "Primera"
$ws.Cells.Item(133, 13).Value = 65
$ws.Cells.Item(133, 14).Value = 3500
$ws.Cells.Item(133, 15).Value = 3500
$ws.Cells.Item(133, 16).Value = 3500
$ws.Cells.Item(133, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(133, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(133, 19).Value = 3500
$ws.Cells.Item(133, 20).Value = 1
